# Adds season "Wins", "Losses" and "Ties" record columns (AD, AE, AF)
# to the player/roster table on Sheet1, matching the format of the
# existing header row and filling the record for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add the three new headers ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font, borders, centered/top aligned)
# from the existing last header cell (AC1) onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (rows 2-41): fill in the team's season record ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 41) { $lastRow = 41 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 72  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 90  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}

$excel.CutCopyMode = 0
